$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 128, shifting rows 128:138 down to 129:139
$ws.Rows.Item(128).Insert()

# Populate the new row 128 with the new weekly record
$ws.Cells.Item(128, 1).Value = 10
$ws.Cells.Item(128, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(128, 3).Value = "La Araucanía"
$ws.Cells.Item(128, 4).Value = 44461
$ws.Cells.Item(128, 4).NumberFormat = $ws.Cells.Item(129, 4).NumberFormat
$ws.Cells.Item(128, 5).Value = 9
$ws.Cells.Item(128, 6).Value = "Fruta"
$ws.Cells.Item(128, 7).Value = 100102
$ws.Cells.Item(128, 8).Value = "Cítricos"
$ws.Cells.Item(128, 9).Value = 100102006
$ws.Cells.Item(128, 10).Value = "Pomelo"
$ws.Cells.Item(128, 11).Value = "Start Ruby"
$ws.Cells.Item(128, 12).Value = "Primera"
$ws.Cells.Item(128, 13).Value = 80
$ws.Cells.Item(128, 14).Value = 12000
$ws.Cells.Item(128, 15).Value = 12000
$ws.Cells.Item(128, 16).Value = 12000
$ws.Cells.Item(128, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(128, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(128, 19).Value = 800
$ws.Cells.Item(128, 20).Value = 15
